$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1603853333333333
$ws.Range("H2").Value = 0.481156
$ws.Range("I2").Value = 0.01032935781992836
$ws.Range("J2").Value = 0.01042870175281933
$ws.Range("M2").Value = 0.8584576666666667
$ws.Range("N2").Value = 2.575373
$ws.Range("O2").Value = 0.02952026538348031
$ws.Range("P2").Value = 0.03028938521394646
$ws.Range("Q2").Value = 0.1376840190208889
$ws.Range("R2").Value = 1.239156171188
$ws.Range("S2").Value = 0.0003049253840852128
$ws.Range("T2").Value = 0.0003158789646725035

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1603853333333333
$ws.Range("H3").Value = 0.481156
$ws.Range("I3").Value = 0.01032935781992836
$ws.Range("J3").Value = 0.01042870175281933
$ws.Range("N3").Value = 75.717583
$ws.Range("O3").Value = 0.8679143348771993
$ws.Range("P3").Value = 0.8905269407406087
$ws.Range("Q3").Value = 4.047996596216445
$ws.Range("R3").Value = 36.431969365948
$ws.Range("S3").Value = 0.00896499772199172
$ws.Range("T3").Value = 0.009287039867834426

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1603853333333333
$ws.Range("H4").Value = 0.481156
$ws.Range("I4").Value = 0.01032935781992836
$ws.Range("J4").Value = 0.01042870175281933
$ws.Range("M4").Value = 0.3580240000000001
$ws.Range("N4").Value = 1.074072
$ws.Range("O4").Value = 0.0123115721415754
$ws.Range("P4").Value = 0.01263233735676886
$ws.Range("Q4").Value = 0.05742179858133334
$ws.Range("R4").Value = 0.516796187232
$ws.Range("S4").Value = 0.000127170633976194
$ws.Range("T4").Value = 0.0001317388787347406

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1603853333333333
$ws.Range("H5").Value = 0.481156
$ws.Range("I5").Value = 0.01032935781992836
$ws.Range("J5").Value = 0.01042870175281933
$ws.Range("M5").Value = 2.2152535
$ws.Range("N5").Value = 4.430507
$ws.Range("O5").Value = 0.07617716487477769
$ws.Range("P5").Value = 0.05210792115009603
$ws.Range("Q5").Value = 0.3552941710153333
$ws.Range("R5").Value = 2.131765026092
$ws.Range("S5").Value = 0.000786861193699257
$ws.Range("T5").Value = 0.0005434179686337783

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1603853333333333
$ws.Range("H6").Value = 0.481156
$ws.Range("I6").Value = 0.01032935781992836
$ws.Range("J6").Value = 0.01042870175281933
$ws.Range("M6").Value = 0.4093533333333334
$ws.Range("N6").Value = 1.22806
$ws.Range("O6").Value = 0.01407666272296744
$ws.Range("P6").Value = 0.01444341553857988
$ws.Range("Q6").Value = 0.06565427081777779
$ws.Range("R6").Value = 0.5908884373600001
$ws.Range("S6").Value = 0.0001454028861759778
$ws.Range("T6").Value = 0.000150626072943886

$ws.Range("I7").Value = 0.9610926076617912
$ws.Range("J7").Value = 0.9703360399430661
$ws.Range("M7").Value = 0.8584576666666667
$ws.Range("N7").Value = 2.575373
$ws.Range("O7").Value = 0.02952026538348031
$ws.Range("P7").Value = 0.03028938521394646
$ws.Range("Q7").Value = 12.81077635037911
$ws.Range("R7").Value = 115.296987153412
$ws.Range("S7").Value = 0.0283717088362772
$ws.Range("T7").Value = 0.02939088210081086

$ws.Range("I8").Value = 0.9610926076617912
$ws.Range("J8").Value = 0.9703360399430661
$ws.Range("N8").Value = 75.717583
$ws.Range("O8").Value = 0.8679143348771993
$ws.Range("P8").Value = 0.8905269407406087
$ws.Range("R8").Value = 3389.803804900652
$ws.Range("S8").Value = 0.8341460513341765
$ws.Range("T8").Value = 0.8641103851408557

$ws.Range("I9").Value = 0.9610926076617912
$ws.Range("J9").Value = 0.9703360399430661
$ws.Range("M9").Value = 0.3580240000000001
$ws.Range("N9").Value = 1.074072
$ws.Range("O9").Value = 0.0123115721415754
$ws.Range("P9").Value = 0.01263233735676886
$ws.Range("Q9").Value = 5.342797403018667
$ws.Range("R9").Value = 48.085176627168
$ws.Range("S9").Value = 0.01183256097396297
$ws.Range("T9").Value = 0.01225761220599196

$ws.Range("I10").Value = 0.9610926076617912
$ws.Range("J10").Value = 0.9703360399430661
$ws.Range("M10").Value = 2.2152535
$ws.Range("N10").Value = 4.430507
$ws.Range("O10").Value = 0.07617716487477769
$ws.Range("P10").Value = 0.05210792115009603
$ws.Range("Q10").Value = 33.05826047088467
$ws.Range("R10").Value = 198.349562825308
$ws.Range("S10").Value = 0.07321331003378229
$ws.Range("T10").Value = 0.05056219385844972

$ws.Range("I11").Value = 0.9610926076617912
$ws.Range("J11").Value = 0.9703360399430661
$ws.Range("M11").Value = 0.4093533333333334
$ws.Range("N11").Value = 1.22806
$ws.Range("O11").Value = 0.01407666272296744
$ws.Range("P11").Value = 0.01444341553857988
$ws.Range("Q11").Value = 6.108785797182223
$ws.Range("R11").Value = 54.97907217464
$ws.Range("S11").Value = 0.01352897648359231
$ws.Range("T11").Value = 0.01401496663695775

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.443735
$ws.Range("H12").Value = 0.88747
$ws.Range("I12").Value = 0.02857803451828042
$ws.Range("J12").Value = 0.01923525830411462
$ws.Range("M12").Value = 0.8584576666666667
$ws.Range("N12").Value = 2.575373
$ws.Range("O12").Value = 0.02952026538348031
$ws.Range("P12").Value = 0.03028938521394646
$ws.Range("Q12").Value = 0.3809277127183333
$ws.Range("R12").Value = 2.28556627631
$ws.Range("S12").Value = 0.0008436311631178988
$ws.Range("T12").Value = 0.0005826241484630902

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.443735
$ws.Range("H13").Value = 0.88747
$ws.Range("I13").Value = 0.02857803451828042
$ws.Range("J13").Value = 0.01923525830411462
$ws.Range("N13").Value = 75.717583
$ws.Range("O13").Value = 0.8679143348771993
$ws.Range("P13").Value = 0.8905269407406087
$ws.Range("Q13").Value = 11.19951389750167
$ws.Range("R13").Value = 67.19708338501
$ws.Range("S13").Value = 0.02480328582103099
$ws.Range("T13").Value = 0.01712951573191859

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 0.443735
$ws.Range("H14").Value = 0.88747
$ws.Range("I14").Value = 0.02857803451828042
$ws.Range("J14").Value = 0.01923525830411462
$ws.Range("M14").Value = 0.3580240000000001
$ws.Range("N14").Value = 1.074072
$ws.Range("O14").Value = 0.0123115721415754
$ws.Range("P14").Value = 0.01263233735676886
$ws.Range("Q14").Value = 0.15886777964
$ws.Range("R14").Value = 0.9532066778400001
$ws.Range("S14").Value = 0.0003518405336362414
$ws.Range("T14").Value = 0.0002429862720421657

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.5
$ws.Range("G15").Value = 0.443735
$ws.Range("H15").Value = 0.88747
$ws.Range("I15").Value = 0.02857803451828042
$ws.Range("J15").Value = 0.01923525830411462
$ws.Range("M15").Value = 2.2152535
$ws.Range("N15").Value = 4.430507
$ws.Range("O15").Value = 0.07617716487477769
$ws.Range("P15").Value = 0.05210792115009603
$ws.Range("Q15").Value = 0.9829855118225
$ws.Range("R15").Value = 3.93194204729
$ws.Range("S15").Value = 0.002176993647296135
$ws.Range("T15").Value = 0.001002309323012535

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.5
$ws.Range("G16").Value = 0.443735
$ws.Range("H16").Value = 0.88747
$ws.Range("I16").Value = 0.02857803451828042
$ws.Range("J16").Value = 0.01923525830411462
$ws.Range("M16").Value = 0.4093533333333334
$ws.Range("N16").Value = 1.22806
$ws.Range("O16").Value = 0.01407666272296744
$ws.Range("P16").Value = 0.01444341553857988
$ws.Range("Q16").Value = 0.1816444013666667
$ws.Range("R16").Value = 1.0898664082
$ws.Range("S16").Value = 0.0004022833531991548
$ws.Range("T16").Value = 0.0002778228286782469
